$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the D column header and add a new E column "PracCorrectPilot"
# with a left/right label for each trial.
$ws.Cells.Item(1, 4).Value = "PracCorrectScan"
$ws.Cells.Item(1, 5).Value = "PracCorrectPilot"

$ws.Cells.Item(2, 5).Value = "left"
$ws.Cells.Item(3, 5).Value = "left"
$ws.Cells.Item(4, 5).Value = "right"
$ws.Cells.Item(5, 5).Value = "left"
$ws.Cells.Item(6, 5).Value = "left"

# Widen the image-path columns (B and C) to fit the longer paths.
$ws.Columns.Item(2).ColumnWidth = 23.8
$ws.Columns.Item(3).ColumnWidth = 23.8

# Move the selection down past the new data, like the author did.
$ws.Range("E7").Select()

# Resize the workbook window, matching the larger view used while editing.
$win = $wb.Windows.Item(1)
$win.Width = 1299
$win.Height = 723
